$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the id sequence: rows 27..101 hold values 26..100 in column A
# (rows 2..26 already hold 1..25).
for ($i = 26; $i -le 100; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
}

# Match the bestFit column widths Excel computed for the header/data text
# (equipementId, equipementLibelle, equipementDescription, equipementVideo,
# typeEquipementId) once the sheet held the full data set.
$ws.Columns.Item(1).ColumnWidth = 11.333333333333334
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 15
